$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.737.13'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.595.69'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.88'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  -2.26%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.40'
$ws.Range("E8").Value = '  -2.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.254'
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0868'
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").Value = '1.822.85'
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("D13").Value = '1.592.86'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("E14").Value = '  -3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.535'
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("D16").Value = '27.731.70'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.57'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.59'
$ws.Range("E18").Value = '  -3.37%  '
$ws.Range("E19").Value = '  -2.25%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  -3.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.78'
$ws.Range("E23").Value = '  -1.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  -4.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.50'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.82'
$ws.Range("E26").Value = '  -1.26%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  -1.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.106'
$ws.Range("E29").Value = '  -4.34%  '
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("E32").Value = '  -4.67%  '
$ws.Range("D33").Value = '1.379.78'
$ws.Range("E33").Value = '  -2.58%  '
$ws.Range("E34").Value = '  -3.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("E35").Value = '  -3.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.973'
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.537'
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.830'
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.23'
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").Value = '1.733.65'
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.01'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("D49").Value = '0.0₆0101'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0968'
$ws.Range("E50").Value = '  -3.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0495'
$ws.Range("E51").Value = '  -1.36%  '
